$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Re-use the date format already applied to the existing date column (A2:A80)
# by copying the formatting (not the value) from the last existing row.
$ws.Range("A80").Copy()
$ws.Range("A81:A83").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New rows of work-hours data appended after the existing last row (row 80)
# Row 81: 2018-12-24, David, "Skype-Konferenz, Besprechung der bevorstehenden Tasks", 1
$ws.Range("A81").Value = 43458
$ws.Range("B81").Value = "David"
$ws.Range("C81").Value = "Skype-Konferenz, Besprechung der bevorstehenden Tasks"
$ws.Range("D81").Value = 1

# Row 82: 2018-12-24, Daniel, "Skype-Konferenz, Besprechung der bevorstehenden Tasks", 1
$ws.Range("A82").Value = 43458
$ws.Range("B82").Value = "Daniel"
$ws.Range("C82").Value = "Skype-Konferenz, Besprechung der bevorstehenden Tasks"
$ws.Range("D82").Value = 1

# Row 83: 2018-12-24, Daniel, "Zeitenprotokoll aktualisiert, Pop-Up für ändern von SessionDate, KV-Translation docx", 3.5
$ws.Range("A83").Value = 43458
$ws.Range("B83").Value = "Daniel"
$ws.Range("C83").Value = "Zeitenprotokoll aktualisiert, Pop-Up für ändern von SessionDate, KV-Translation docx"
$ws.Range("D83").Value = 3.5

# Update selection to match final state of the file
$ws.Range("F81").Select()
